# Insert a new data row at row 51 (pushes existing rows 51..146 down to 52..147)
# and populate it with the new weekly price observation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("51:51").Insert()

$ws.Range("A51").Value = 4
$ws.Range("B51").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C51").Value = "Los Lagos"
$ws.Range("D51").Value = 44533
$ws.Range("E51").Value = 10
$ws.Range("F51").Value = 100112028
$ws.Range("G51").Value = "Sandia"
$ws.Range("H51").Value = "Sin especificar"
$ws.Range("I51").Value = "Primera"
$ws.Range("J51").Value = 700
$ws.Range("K51").Value = 900
$ws.Range("L51").Value = 1000
$ws.Range("M51").Value = 950
$ws.Range("N51").Value = "$/kilo (volumen en unidades)"
$ws.Range("O51").Value = "Perú"
$ws.Range("P51").Value = 950
$ws.Range("Q51").Value = 1
$ws.Range("R51").Value = "Hortaliza"
